$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the first "SourceCode" paragraph: insert a brand new paragraph in
#    front of it (style "FirstParagraph") carrying the new explanatory text,
#    and fix up the first line of the remaining SourceCode paragraph
#    ("Summary 2" -> "Summary 0", and a double space collapsed to one).
# ---------------------------------------------------------------------------

$sourcePara = $d.Paragraphs.Item(2)
$insertionPoint = $sourcePara.Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item(2)
$newStart = $newPara.Range.Start
$fullText = "When su_label(cols) and cat_col are both specified the column containing the catagory names comes after the column with the summary labels."
$typeRng = $newPara.Range.Duplicate
$typeRng.Collapse(1)
$typeRng.InsertAfter($fullText)
$newEnd = $newStart + $fullText.Length

# Newly inserted runs pick up whatever character style happens to be "active"
# (VerbatimChar here) -- explicitly reset the whole new paragraph to the
# plain/default character style first ...
$d.Range($newStart, $newEnd).Style = "Default Paragraph Font"

# ... then re-apply the VerbatimChar style to the two code tokens.
$tok1 = $d.Range($newStart, $newEnd)
$tok1.Find.Execute("su_label(cols)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tok1.Style = "Verbatim Char"

$tok2 = $d.Range($tok1.End, $newEnd)
$tok2.Find.Execute("cat_col", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tok2.Style = "Verbatim Char"

# Finally apply the paragraph style to the new paragraph.
$newPara.Style = "First Paragraph"

# Fix up the text of the (now third) SourceCode paragraph's first line.
$summaryRng = $d.Content
$summaryRng.Find.Execute("(""Summary 1"")  (""Summary 2"")", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$summaryRng.Text = "(""Summary 1"") (""Summary 0"")"

# ---------------------------------------------------------------------------
# 2) Remove the "gap(2)" option from the "ethnicity" example line.
# ---------------------------------------------------------------------------

$gapRng = $d.Content
$gapRng.Find.Execute("cat_levels(4 3 2 1 0)  gap(2)   su_label(col)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gapRng.Text = "cat_levels(4 3 2 1 0)   su_label(col)"
